$d = $word.ActiveDocument

# Locate the target text in the "Location/Condition/Action" alternative-flow
# table: "If the balance is not enough" -> "If the mandatory field is left
# blank", split across two runs ("If the " / "mandatory field is left blank").
$rng = $d.Content
$found = $rng.Find.Execute("If the balance is not enough", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'If the balance is not enough'"
}

$start = $rng.Start
$end = $rng.End

# Replace the whole found range with the final text first (keeps it as a
# single run with the original run formatting).
$full = $d.Range($start, $end)
$full.Text = "If the mandatory field is left blank"

$splitAt = $start + "If the ".Length

# Force a run split at the boundary between "If the " and "mandatory field
# is left blank" by toggling character formatting on the second part and
# then reverting it - adjacent runs with identical formatting are merged,
# but this round-trip leaves a genuine run boundary in place.
$second = $d.Range($splitAt, $start + "If the mandatory field is left blank".Length)
$second.Bold = 1
$second.Bold = 0
